$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the formatting of the other
# header cells (e.g. G1) by copying G1's formatting into H1 first.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the "Save" value for the data row
$ws.Range("H2").Value = 1
